$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the coinranking.com crypto snapshot refresh.
# Numeric-looking text values (prices without thousands-dot formatting,
# e.g. "583.44") need an explicit Text format while assigning, otherwise
# Excel auto-converts the string to a number; the style is then reset
# back to Normal so no stray number-format sticks to the cell.

$ws.Range("D2").Value = "63.119.28"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.603.20"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "3.070.57"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "63.050.53"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "2.604.73"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").Value = "2.725.36"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.53%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.98%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "464.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.406"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.01%  "
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "159.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.32%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.640"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0547"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0975"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("E51").Value = "  -0.12%  "
